$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 3.55
$ws.Range("L2").Value = 1.34
$ws.Range("V2").Value = 1.54
$ws.Range("Y2").Value = 16
$ws.Range("AL2").Value = 44

# Row 3
$ws.Range("P3").Value = 1.89

# Row 4
$ws.Range("Q4").Value = 1.7
$ws.Range("T4").Value = 1.52

# Row 6
$ws.Range("Q6").Value = 1.81
$ws.Range("T6").Value = 1.67

# Row 7
$ws.Range("T7").Value = 1.69
$ws.Range("Y7").Value = 1000
$ws.Range("AC7").Value = 13
